$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cell corrections (NBA stats off-by-one-day fix)
$ws.Range("AH2").Value = 7
$ws.Range("AR2").Value = 27
$ws.Range("AO3").Value = 21
$ws.Range("AD4").Value = 7
$ws.Range("AP6").Value = 21
$ws.Range("D8").Value = 63
$ws.Range("F8").Value = 33
$ws.Range("G8").Value = 0.476
$ws.Range("O8").Value = 16.9
$ws.Range("P8").Value = 21.4
$ws.Range("Q8").Value = 0.792
$ws.Range("R8").Value = 9.4
$ws.Range("T8").Value = 42.3
$ws.Range("V8").Value = 14.2
$ws.Range("X8").Value = 5.4
$ws.Range("AA8").Value = 19.4
$ws.Range("AB8").Value = 101.4
$ws.Range("AD8").Value = 23
$ws.Range("AI8").Value = 4
$ws.Range("AN8").Value = 8
$ws.Range("AO8").Value = 15
$ws.Range("AP8").Value = 20
$ws.Range("AR8").Value = 26
$ws.Range("AT8").Value = 15
$ws.Range("AH9").Value = 11
$ws.Range("AD12").Value = 7
$ws.Range("AI12").Value = 5
$ws.Range("AS13").Value = 5
$ws.Range("AO14").Value = 19
$ws.Range("AD16").Value = 23
$ws.Range("AJ16").Value = 14
$ws.Range("AN16").Value = 23
$ws.Range("AD17").Value = 23
$ws.Range("AE17").Value = 1
$ws.Range("AD18").Value = 23
$ws.Range("AD19").Value = 29
$ws.Range("AD20").Value = 7
$ws.Range("AN20").Value = 7
$ws.Range("D21").Value = 62
$ws.Range("F21").Value = 24
$ws.Range("G21").Value = 0.613
$ws.Range("N21").Value = 0.368
$ws.Range("O21").Value = 16.6
$ws.Range("P21").Value = 21.8
$ws.Range("Q21").Value = 0.76
$ws.Range("W21").Value = 8.199999999999999
$ws.Range("Z21").Value = 19.5
$ws.Range("AA21").Value = 19.5
$ws.Range("AB21").Value = 99.3
$ws.Range("AC21").Value = 3.2
$ws.Range("AD21").Value = 29
$ws.Range("AF21").Value = 7
$ws.Range("AN21").Value = 9
$ws.Range("AO21").Value = 18
$ws.Range("AS21").Value = 21
$ws.Range("AW21").Value = 12
$ws.Range("AD22").Value = 7
$ws.Range("AD23").Value = 7
$ws.Range("AT23").Value = 16
$ws.Range("AX24").Value = 21
$ws.Range("AD25").Value = 7
$ws.Range("D26").Value = 63
$ws.Range("E26").Value = 29
$ws.Range("G26").Value = 0.46
$ws.Range("H26").Value = 48.6
$ws.Range("I26").Value = 36.7
$ws.Range("J26").Value = 82
$ws.Range("K26").Value = 0.448
$ws.Range("L26").Value = 8.1
$ws.Range("Q26").Value = 0.778
$ws.Range("S26").Value = 30.3
$ws.Range("T26").Value = 41.4
$ws.Range("U26").Value = 21.5
$ws.Range("Z26").Value = 18.8
$ws.Range("AB26").Value = 98
$ws.Range("AC26").Value = -1.8
$ws.Range("AD26").Value = 23
$ws.Range("AE26").Value = 19
$ws.Range("AF26").Value = 19
$ws.Range("AG26").Value = 19
$ws.Range("AH26").Value = 5
$ws.Range("AJ26").Value = 13
$ws.Range("AN26").Value = 24
$ws.Range("AO26").Value = 20
$ws.Range("AX26").Value = 22
$ws.Range("AZ26").Value = 7
$ws.Range("BC26").Value = 20
$ws.Range("AN27").Value = 10
$ws.Range("D28").Value = 65
$ws.Range("E28").Value = 49
$ws.Range("G28").Value = 0.754
$ws.Range("I28").Value = 39.6
$ws.Range("J28").Value = 81.2
$ws.Range("L28").Value = 8.4
$ws.Range("M28").Value = 21.9
$ws.Range("O28").Value = 16.9
$ws.Range("S28").Value = 33
$ws.Range("T28").Value = 40.9
$ws.Range("U28").Value = 25.1
$ws.Range("Y28").Value = 4.6
$ws.Range("Z28").Value = 17.6
$ws.Range("AA28").Value = 19
$ws.Range("AB28").Value = 104.4
$ws.Range("AC28").Value = 8.1
$ws.Range("AD28").Value = 7
$ws.Range("AH28").Value = 8
$ws.Range("AO28").Value = 16
$ws.Range("AY28").Value = 9
$ws.Range("AD29").Value = 7
$ws.Range("AX29").Value = 20
$ws.Range("BC29").Value = 19
$ws.Range("AD30").Value = 7
$ws.Range("AW30").Value = 11
$ws.Range("AD31").Value = 23
$ws.Range("AY31").Value = 10

# Fix Date column format for all data rows (BF2:BF31): "3-14-2012-13" -> "2013-03-14"
# Force text format first so Excel does not reinterpret the string as a date serial.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.NumberFormat = "@"
    $cell.Value = "2013-03-14"
}
